# Apply the commit: add two new rows ("Holden" and "Rizzie Spiral") after
# the "Spiral5" row, rename "Thomas Hex" -> "Matthies Hex", and rerun the
# simulation (new data values) for all rows below the insertion point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows directly below row 3 ("Spiral5"), shifting all
#     subsequent rows (old row 4 .. old row 29) down by two rows.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Copy formatting (bold/border/centered style) from the row-3 label cell
# onto the two freshly inserted label cells so they match the rest of the
# table (Insert() alone leaves them with a slightly different auto style).
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A3").Copy($ws.Range("A5"))

# --- Row 4: new series "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.8115741179572922
$ws.Range("D4").Value = 0.9482788694071127
$ws.Range("E4").Value = 1.093708900492132
$ws.Range("F4").Value = 0.8115741179572922
$ws.Range("G4").Value = 1.048913678219783
$ws.Range("H4").Value = 0.9482788694071127
$ws.Range("I4").Value = 0.8493825035583669
$ws.Range("J4").Value = 1.403113301400745
$ws.Range("K4").Value = 0.9482788694071127
$ws.Range("L4").Value = 1.093708900492132
$ws.Range("M4").Value = 0.9526415092247121
$ws.Range("N4").Value = 0.9526415092247121
$ws.Range("O4").Value = 0.9182218406692636
$ws.Range("P4").Value = 0.9511872959521789
$ws.Range("Q4").Value = 0.9511872959521789
$ws.Range("R4").Value = 0.9504601893159124
$ws.Range("S4").Value = 0.9504601893159124
$ws.Range("T4").Value = 1.025828561839239

# --- Row 5: new series "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.285660442863944
$ws.Range("D5").Value = 0.7837792567207699
$ws.Range("E5").Value = 1.120037499514067
$ws.Range("F5").Value = 1.285660442863944
$ws.Range("G5").Value = 0.9926907483397743
$ws.Range("H5").Value = 0.7837792567207699
$ws.Range("I5").Value = 0.8503098473864319
$ws.Range("J5").Value = 1.687968925597128
$ws.Range("K5").Value = 0.7837792567207699
$ws.Range("L5").Value = 1.120037499514067
$ws.Range("M5").Value = 1.202848971189006
$ws.Range("N5").Value = 1.202848971189006
$ws.Range("O5").Value = 1.085335929921481
$ws.Range("P5").Value = 1.06315906636626
$ws.Range("Q5").Value = 1.06315906636626
$ws.Range("R5").Value = 0.9933141139548878
$ws.Range("S5").Value = 0.9933141139548878
$ws.Range("T5").Value = 1.120074453403686

# --- Rename the "Thomas Hex" series to "Matthies Hex" (old row 9, now row
#     11 after the two-row insertion above).
$ws.Range("B11").Value = "Matthies Hex"
